# response_to_reviewer_R1.docx - apply reviewer-response edits
# Strategy: for every paragraph we touch, clear its run text to "" first
# (this drops the old <w:r>/<w:rPr> entirely, leaving a clean <w:p><w:pPr>.../></w:p>)
# and then re-type the new text into the now-empty paragraph so the new
# run starts from a blank slate instead of inheriting stray formatting.
# We work from the bottom of the document upwards so paragraph indices
# for edits still to come are not disturbed by earlier edits.

$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $text, $bold) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $rNoMark = $d.Range($r.Start, $r.End - 1)
    $rNoMark.Text = ""

    $p2 = $d.Paragraphs.Item($paraIndex)
    $r2 = $p2.Range
    $r2NoMark = $d.Range($r2.Start, $r2.End - 1)
    $r2NoMark.Text = $text
    if ($bold) {
        $r2NoMark.Bold = 1
    }
    return $r2NoMark
}

# ---------------------------------------------------------------------
# Hunk 5: paragraph 40 (crown exposure RESPONSE) + new paragraph after it
# ---------------------------------------------------------------------
$p40 = $d.Paragraphs.Item(40)
$r40 = $p40.Range
$r40NoMark = $d.Range($r40.Start, $r40.End - 1)
$r40NoMark.Text = ""

$p40b = $d.Paragraphs.Item(40)
$r40b = $p40b.Range
$r40bNoMark = $d.Range($r40b.Start, $r40b.End - 1)
$r40bNoMark.Text = "RESPONSE: Since all crown exposure categories had negative growth sensitivities in 2015 and 2020, crown exposure was not associated with large differences in the magnitude or direction of sensitivity. In other words, more exposed trees had more or less similar growth declines compared to less exposed/understory trees. We have now clarified this in the text which now reads:"
$r40bNoMark.Bold = 1

# Insert a brand-new BodyText paragraph right after paragraph 40 containing
# the quoted sentence (three runs: opening quote / sentence / closing quote)
$p40c = $d.Paragraphs.Item(40)
$p40c.Range.InsertParagraphAfter()

$pNew = $d.Paragraphs.Item(41)
$rNew = $pNew.Range
$rNewNoMark = $d.Range($rNew.Start, $rNew.End - 1)
$rNewNoMark.Text = "“Predicted effects of CII on sensitivity decreased monotonically, with a decrease to negative sensitivities in category 4 and 5 in 2010, but all CII categories in 2015 and 2020 has similar predictions, showing that exposure was not associated with altered sensitivity.”"

Write-Output "hunk5 done"

# ---------------------------------------------------------------------
# Hunk 4: paragraph 31 "RESPONSE: The three climatically distinct droughts"
# ---------------------------------------------------------------------
Set-ParaText 31 "RESPONSE: The three climatically distinct droughts are summarised in Table 1. In summary, 2010 was a moderate drought that peaked in the dry season, 2015 was a severe drought that peaked in the wet season and 2020 was a moderate drought that peaked in the transition from dry to wet season. After moving the drought descriptions to the methods section, we have added this summary to set context for the results section." $true | Out-Null
Write-Output "hunk4 done"

# ---------------------------------------------------------------------
# Hunk 3: paragraph 26 - supplementary table / IUCN categories response
# ---------------------------------------------------------------------
$p26 = $d.Paragraphs.Item(26)
$r26 = $p26.Range
$r26NoMark = $d.Range($r26.Start, $r26.End - 1)
$r26NoMark.Text = ""

$p26b = $d.Paragraphs.Item(26)
$r26b = $p26b.Range
$r26bNoMark = $d.Range($r26b.Start, $r26b.End - 1)
$r26bNoMark.Text = "RESPONSE: Thanks for this suggestion. We have now included a supplementary table (Table S1) that details the response of specific species to specific drought events. This summary highlights our main message that the species responses are dissimilar across drought events. This table also includes available information on the IUCN categories of the species analysed (accessed in September 2025) where information was avaialable. Two of the species analysed are “Endangered”, including "
$r26bNoMark.Bold = 1

$afterFirst = $r26bNoMark.End
$italicRange = $d.Range($afterFirst, $afterFirst)
$italicRange.InsertAfter("Afzelia xylocarpa")
$italicRange2 = $d.Range($afterFirst, $afterFirst + ("Afzelia xylocarpa".Length))
$italicRange2.Bold = 1
$italicRange2.Italic = 1

$afterItalic = $italicRange2.End
$tailRange = $d.Range($afterItalic, $afterItalic)
$tailText = ", the species that was discussed. There are also “Vulnerable” and “Data Deficient” species in the list, some of which have consistent negative drought sensitivities. This table would therefore be useful for forest managers and conservation professionals and to guide species-specific research directions."
$tailRange.InsertAfter($tailText)
$tailRange2 = $d.Range($afterItalic, $afterItalic + $tailText.Length)
$tailRange2.Bold = 1

Write-Output "hunk3 done"

# ---------------------------------------------------------------------
# Hunk 2: paragraph 20 - empty RESPONSE: after "Predictions (i-iii)" comment
# ---------------------------------------------------------------------
Set-ParaText 20 "RESPONSE: Thanks for this suggestion. We have now revised the Results section to remove some redundant statistical details that are reported in figures. We have also reorganised some text so that the statements that correspond to each prediction are highlighted at the start of each section. We hope that the results and their alignment with the predictions are clearer now." $true | Out-Null
Write-Output "hunk2 done"

# ---------------------------------------------------------------------
# Hunk 1: paragraph 18 - empty RESPONSE: after "acronyms" comment
# ---------------------------------------------------------------------
Set-ParaText 18 "RESPONSE: We have now carefully reviewed the text and rewritten the description of acronyms (TWI and CII) in the results section to make it clearer for readers." $true | Out-Null
Write-Output "hunk1 done"

Write-Output "all done"
